# Change date format tokens from strftime-less (yyyymmdd...) style to
# Python-strftime-style (%Y%m%d...) in the "Parameters" column (F) of the
# CDAVariables_short sheet, per commit:
#   "Change date format from yyyymmddhhmmss to %Y%m%d%H%M%S"
#
# Several of the target cells end up as rich text (multiple runs) because,
# in the real edit, only the date-format token substring was retyped while
# the rest of the text (and its original run formatting) was left alone.
# We reproduce that by using Range.Characters(start,len) to restyle just
# the replaced token, which forces the engine to split the string into
# runs the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "date_creation" parameters (F3) --------------------------------
# Plain text replacement only -- no run split in the target diff.
$ws.Range("F3").Value2 = "date_format=%Y%m%d%H%M%S;start_date=20210101;end_date=20240730"

# --- Row 13: "author_ts" parameters (F13) ----------------------------------
# Ends up with the same visible text as F3 used to (and now does again),
# but as its own distinct rich-text string (3 runs).
$cell = $ws.Range("F13")
$cell.Value2 = "date_format=%Y%m%d%H%M%S;start_date=20210101;end_date=20240730"
$mid = $cell.Characters(13, 12)
$mid.Font.Name = "Calibri"
$mid.Font.Size = 11

# --- Row 12: "birth_date_ts" parameters (F12) ------------------------------
$cell = $ws.Range("F12")
$cell.Value2 = "date_format=%Y%m%d;start_date=19200101;end_date=20200730"
$mid = $cell.Characters(13, 6)
$mid.Font.Name = "Calibri"
$mid.Font.Size = 11

# --- Row 26: "discharge_ts" parameters (F26) -------------------------------
$cell = $ws.Range("F26")
$cell.Value2 = "date_format=%Y%m%d%H%M%S"
$mid = $cell.Characters(13, 12)
$mid.Font.Name = "Calibri"
$mid.Font.Size = 11

# --- Row 36: "complaint_start" parameters (F36) -----------------------------
$cell = $ws.Range("F36")
$cell.Value2 = "date_format=%Y%m%d;start_date=20210101;end_date=20240730"
$mid = $cell.Characters(13, 6)
$mid.Font.Name = "Calibri"
$mid.Font.Size = 11

# --- Row 41: "triage_ts_start" parameters (F41) -----------------------------
$cell = $ws.Range("F41")
$cell.Value2 = "date_format=%Y%m%d%H%M"
$mid = $cell.Characters(13, 10)
$mid.Font.Name = "Calibri"
$mid.Font.Size = 11

# --- Row 50: "diagnostics_ts" parameters (F50) ------------------------------
$cell = $ws.Range("F50")
$cell.Value2 = "date_format=%Y%m%d%H%M;start_date=20210101;end_date=20240730"
$mid = $cell.Characters(13, 10)
$mid.Font.Name = "Calibri"
$mid.Font.Size = 11

# --- Selection moved from F20 to F19 in the saved view ----------------------
$ws.Range("F19").Select()
